# update scripts wuth new tpm
# Refresh NATMI Fn1-Col13a1 LR-pair metrics (columns G-T, rows 2-13) with
# recomputed values based on the updated TPM expression data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6763629999999999
$ws.Range("N2").Value = 2.029089
$ws.Range("O2").Value = 0.6127318215515719
$ws.Range("P2").Value = 0.6127318215515719
$ws.Range("Q2").Value = 19.765657606898
$ws.Range("R2").Value = 177.890918462082
$ws.Range("S2").Value = 0.01035769494712173
$ws.Range("T2").Value = 0.01035769494712173

$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("M3").Value = 0.05377866666666667
$ws.Range("O3").Value = 0.04871925339984812
$ws.Range("P3").Value = 0.04871925339984811
$ws.Range("Q3").Value = 1.571597961285333
$ws.Range("R3").Value = 14.144381651568
$ws.Range("S3").Value = 0.0008235563210824325
$ws.Range("T3").Value = 0.0008235563210824324

$ws.Range("G4").Value = 29.223446
$ws.Range("H4").Value = 87.670338
$ws.Range("I4").Value = 0.0169041244192178
$ws.Range("J4").Value = 0.0169041244192178
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3737066666666666
$ws.Range("N4").Value = 1.12112
$ws.Range("O4").Value = 0.3385489250485801
$ws.Range("P4").Value = 0.33854892504858
$ws.Range("Q4").Value = 10.92099659317333
$ws.Range("R4").Value = 98.28896933855999
$ws.Range("S4").Value = 0.00572287315101364
$ws.Range("T4").Value = 0.005722873151013639

$ws.Range("I5").Value = 0.9471112884046843
$ws.Range("J5").Value = 0.9471112884046842
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6763629999999999
$ws.Range("N5").Value = 2.029089
$ws.Range("O5").Value = 0.6127318215515719
$ws.Range("P5").Value = 0.6127318215515719
$ws.Range("Q5").Value = 1107.438455726963
$ws.Range("R5").Value = 9966.94610154267
$ws.Range("S5").Value = 0.5803252249562584
$ws.Range("T5").Value = 0.5803252249562583

$ws.Range("I6").Value = 0.9471112884046843
$ws.Range("J6").Value = 0.9471112884046842
$ws.Range("M6").Value = 0.05377866666666667
$ws.Range("O6").Value = 0.04871925339984812
$ws.Range("P6").Value = 0.04871925339984811
$ws.Range("Q6").Value = 88.05414188000888
$ws.Range("R6").Value = 792.48727692008
$ws.Range("S6").Value = 0.04614255485764444
$ws.Range("T6").Value = 0.04614255485764443

$ws.Range("I7").Value = 0.9471112884046843
$ws.Range("J7").Value = 0.9471112884046842
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3737066666666666
$ws.Range("N7").Value = 1.12112
$ws.Range("O7").Value = 0.3385489250485801
$ws.Range("P7").Value = 0.33854892504858
$ws.Range("Q7").Value = 611.8861230259555
$ws.Range("R7").Value = 5506.975107233599
$ws.Range("S7").Value = 0.3206435085907816
$ws.Range("T7").Value = 0.3206435085907815

$ws.Range("G8").Value = 37.39212666666667
$ws.Range("H8").Value = 112.17638
$ws.Range("I8").Value = 0.02162924801792661
$ws.Range("J8").Value = 0.0216292480179266
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.6763629999999999
$ws.Range("N8").Value = 2.029089
$ws.Range("O8").Value = 0.6127318215515719
$ws.Range("P8").Value = 0.6127318215515719
$ws.Range("Q8").Value = 25.29065096864667
$ws.Range("R8").Value = 227.61585871782
$ws.Range("S8").Value = 0.0132529285368149
$ws.Range("T8").Value = 0.0132529285368149

$ws.Range("G9").Value = 37.39212666666667
$ws.Range("H9").Value = 112.17638
$ws.Range("I9").Value = 0.02162924801792661
$ws.Range("J9").Value = 0.0216292480179266
$ws.Range("M9").Value = 0.05377866666666667
$ws.Range("O9").Value = 0.04871925339984812
$ws.Range("P9").Value = 0.04871925339984811
$ws.Range("Q9").Value = 2.010898715964445
$ws.Range("R9").Value = 18.09808844368
$ws.Range("S9").Value = 0.001053760815033529
$ws.Range("T9").Value = 0.001053760815033529

$ws.Range("G10").Value = 37.39212666666667
$ws.Range("H10").Value = 112.17638
$ws.Range("I10").Value = 0.02162924801792661
$ws.Range("J10").Value = 0.0216292480179266
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3737066666666666
$ws.Range("N10").Value = 1.12112
$ws.Range("O10").Value = 0.3385489250485801
$ws.Range("P10").Value = 0.33854892504858
$ws.Range("Q10").Value = 13.97368701617778
$ws.Range("R10").Value = 125.7631831456
$ws.Range("S10").Value = 0.007322558666078185
$ws.Range("T10").Value = 0.007322558666078182

$ws.Range("G11").Value = 24.817167
$ws.Range("H11").Value = 74.45150100000001
$ws.Range("I11").Value = 0.01435533915817136
$ws.Range("J11").Value = 0.01435533915817136
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.6763629999999999
$ws.Range("N11").Value = 2.029089
$ws.Range("O11").Value = 0.6127318215515719
$ws.Range("P11").Value = 0.6127318215515719
$ws.Range("Q11").Value = 16.785413523621
$ws.Range("R11").Value = 151.068721712589
$ws.Range("S11").Value = 0.008795973111376948
$ws.Range("T11").Value = 0.008795973111376946

$ws.Range("G12").Value = 24.817167
$ws.Range("H12").Value = 74.45150100000001
$ws.Range("I12").Value = 0.01435533915817136
$ws.Range("J12").Value = 0.01435533915817136
$ws.Range("M12").Value = 0.05377866666666667
$ws.Range("O12").Value = 0.04871925339984812
$ws.Range("P12").Value = 0.04871925339984811
$ws.Range("Q12").Value = 1.334634151704
$ws.Range("R12").Value = 12.011707365336
$ws.Range("S12").Value = 0.0006993814060877129
$ws.Range("T12").Value = 0.0006993814060877127

$ws.Range("G13").Value = 24.817167
$ws.Range("H13").Value = 74.45150100000001
$ws.Range("I13").Value = 0.01435533915817136
$ws.Range("J13").Value = 0.01435533915817136
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.3737066666666666
$ws.Range("N13").Value = 1.12112
$ws.Range("O13").Value = 0.3385489250485801
$ws.Range("P13").Value = 0.33854892504858
$ws.Range("Q13").Value = 9.274340755679999
$ws.Range("R13").Value = 83.46906680112001
$ws.Range("S13").Value = 0.004859984640706703
$ws.Range("T13").Value = 0.004859984640706702
